$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.583.50"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.038.75"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'202.06"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'632.65"
$ws.Range("E6").Value = "  +6.06%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  +6.20%  "
$ws.Range("D10").Value = "3.033.98"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").Value = "'0.438"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "'5.19"
$ws.Range("E13").Value = "  +6.69%  "
$ws.Range("D14").Value = "3.591.08"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "'29.62"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").Value = "76.493.66"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("D18").Value = "3.014.33"
$ws.Range("E18").Value = "  +3.74%  "
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "'9.05"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").Value = "'376.98"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "'73.86"
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("D25").Value = "3.184.07"
$ws.Range("E25").Value = "  +4.60%  "
$ws.Range("E26").Value = "  +5.29%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").Value = "'0.0000114"
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'8.34"
$ws.Range("E31").Value = "  +8.29%  "
$ws.Range("D32").Value = "'1.42"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'518.07"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("E34").Value = "  +9.01%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'20.91"
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("D37").Value = "'163.00"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'0.386"
$ws.Range("E38").Value = "  +10.97%  "
$ws.Range("D39").Value = "'20.02"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  +5.42%  "
$ws.Range("D41").Value = "'187.84"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "'5.18"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("D45").Value = "'42.22"
$ws.Range("E45").Value = "  +5.10%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.68"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "'1.26"
$ws.Range("E47").Value = "  +6.34%  "
$ws.Range("E48").Value = "  +6.09%  "
$ws.Range("D49").Value = "'0.728"
$ws.Range("E49").Value = "  +11.19%  "
$ws.Range("D50").Value = "'0.610"
$ws.Range("E50").Value = "  +6.92%  "
$ws.Range("D51").Value = "'3.92"
$ws.Range("E51").Value = "  +5.51%  "
